$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Transportation"
$ws.Range("B2").Value = "smed"
$ws.Range("C2").Value = "'2023-03-20"
$ws.Range("D2").Value = "'453.0"
$ws.Range("E2").Value = "Checkings"

$ws.Range("A3").Value = "Transportation"
$ws.Range("B3").Value = "buss"
$ws.Range("C3").Value = "'2023-03-20"
$ws.Range("D3").Value = "'700.0"
$ws.Range("E3").Value = "Checkings"

$ws.Range("A4").Value = "Transportation"
$ws.Range("B4").Value = "asdf"
$ws.Range("C4").Value = "'2023-03-22"
$ws.Range("D4").Value = "'1000.0"
$ws.Range("E4").Value = "Checkings"

$ws.Range("A2:E4").ClearFormats()
